# The commit swaps the two theme parts in the package: the slide-design
# theme (ppt/theme/theme2.xml, currently the "Integral" palette, used by
# slideMaster1 / the presentation itself) is recoloured to the "Office
# Theme" palette that used to live in ppt/theme/theme1.xml (the notes
# master's theme). Font scheme and format scheme are identical between
# the two theme parts, so the only observable difference is the 12-slot
# theme colour scheme - swap those via the Design's ThemeColorScheme.

function ConvertTo-BGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours (previously theme1.xml).
$officeTheme = [ordered]@{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$scheme = $design.SlideMaster.Theme.ThemeColorScheme

foreach ($idx in $officeTheme.Keys) {
    $scheme.Item($idx).RGB = ConvertTo-BGR $officeTheme[$idx]
}
